$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14 (existing row 14 and below shift down to 15+)
$ws.Rows.Item(14).Insert()

# Populate I15 first (the "Use 2019 data" note for the existing Traffic Density /
# tmas-2020 row that shifted from 14 -> 15) so shared-string insertion order matches.
$ws.Range("I15").Value = "Use 2019 data"

# New row 14: a second "Traffic Density" entry (FHWA HPMS shapefiles)
$ws.Range("A14").Value = "Traffic Density"
$ws.Range("B14").Value = "CONUS, AK, HI"
$ws.Range("C14").Value = "Y"
$ws.Range("D14").Value = ".shp"
$ws.Range("G14").Value = "FHA"
$ws.Range("H14").Value = "https://www.fhwa.dot.gov/policyinformation/hpms/shapefiles.cfm"
$ws.Range("I14").Value = "Methods from EJScreen and CO Enviroscreen: https://www.epa.gov/sites/default/files/2021-04/documents/ejscreen_technical_document.pdf "

# Turn H15 (tmas-2020 link) into a real hyperlink, styled like the other links
$ws.Hyperlinks.Add($ws.Range("H15"), "https://catalog.data.gov/dataset/tmas-2020")
$ws.Range("H15").Style = "Hyperlink"

# Update the active selection to match the author's final cursor position
$ws.Range("G15").Select() | Out-Null

Write-Host "done"
